$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "78+16="
$t.Cell(1,2).Range.Text = "14+4="
$t.Cell(1,3).Range.Text = "35-11="
$t.Cell(1,4).Range.Text = "1+5="
$t.Cell(1,5).Range.Text = "30+2="
$t.Cell(2,1).Range.Text = "67-11="
$t.Cell(2,2).Range.Text = "33-5="
$t.Cell(2,3).Range.Text = "30+0="
$t.Cell(2,4).Range.Text = "0+93="
$t.Cell(2,5).Range.Text = "45+35="
$t.Cell(3,1).Range.Text = "90-86="
$t.Cell(3,2).Range.Text = "48+10="
$t.Cell(3,3).Range.Text = "77+2="
$t.Cell(3,4).Range.Text = "48+51="
$t.Cell(3,5).Range.Text = "3+59="
$t.Cell(4,1).Range.Text = "26+58="
$t.Cell(4,2).Range.Text = "44-40="
$t.Cell(4,3).Range.Text = "41-8="
$t.Cell(4,4).Range.Text = "25+3="
$t.Cell(4,5).Range.Text = "29+66="
$t.Cell(5,1).Range.Text = "63+29="
$t.Cell(5,2).Range.Text = "95-16="
$t.Cell(5,3).Range.Text = "59-56="
$t.Cell(5,4).Range.Text = "3+9="
$t.Cell(5,5).Range.Text = "3+88="
$t.Cell(6,1).Range.Text = "11+80="
$t.Cell(6,2).Range.Text = "3+12="
$t.Cell(6,3).Range.Text = "31-19="
$t.Cell(6,4).Range.Text = "20-0="
$t.Cell(6,5).Range.Text = "13+4="
$t.Cell(7,1).Range.Text = "93-25="
$t.Cell(7,2).Range.Text = "32+16="
$t.Cell(7,3).Range.Text = "88-62="
$t.Cell(7,4).Range.Text = "83-32="
$t.Cell(7,5).Range.Text = "17-11="
$t.Cell(8,1).Range.Text = "80-62="
$t.Cell(8,2).Range.Text = "14+67="
$t.Cell(8,3).Range.Text = "45+6="
$t.Cell(8,4).Range.Text = "55-27="
$t.Cell(8,5).Range.Text = "75-51="
$t.Cell(9,1).Range.Text = "91-76="
$t.Cell(9,2).Range.Text = "74-9="
$t.Cell(9,3).Range.Text = "68-15="
$t.Cell(9,4).Range.Text = "36-32="
$t.Cell(9,5).Range.Text = "10+49="
$t.Cell(10,1).Range.Text = "28+14="
$t.Cell(10,2).Range.Text = "82-15="
$t.Cell(10,3).Range.Text = "50-33="
$t.Cell(10,4).Range.Text = "66+4="
$t.Cell(10,5).Range.Text = "34-10="
$t.Cell(11,1).Range.Text = "19-3="
$t.Cell(11,2).Range.Text = "16+45="
$t.Cell(11,3).Range.Text = "56-9="
$t.Cell(11,4).Range.Text = "51-35="
$t.Cell(11,5).Range.Text = "70-67="
$t.Cell(12,1).Range.Text = "51-8="
$t.Cell(12,2).Range.Text = "11+40="
$t.Cell(12,3).Range.Text = "44+48="
$t.Cell(12,4).Range.Text = "86-32="
$t.Cell(12,5).Range.Text = "52-17="
$t.Cell(13,1).Range.Text = "20+43="
$t.Cell(13,2).Range.Text = "62+15="
$t.Cell(13,3).Range.Text = "61-44="
$t.Cell(13,4).Range.Text = "21+28="
$t.Cell(13,5).Range.Text = "6+76="
$t.Cell(14,1).Range.Text = "52-13="
$t.Cell(14,2).Range.Text = "26-5="
$t.Cell(14,3).Range.Text = "21+35="
$t.Cell(14,4).Range.Text = "90-6="
$t.Cell(14,5).Range.Text = "55-13="
$t.Cell(15,1).Range.Text = "33+16="
$t.Cell(15,2).Range.Text = "68+17="
$t.Cell(15,3).Range.Text = "8+73="
$t.Cell(15,4).Range.Text = "5-2="
$t.Cell(15,5).Range.Text = "37+43="
$t.Cell(16,1).Range.Text = "72-5="
$t.Cell(16,2).Range.Text = "64-43="
$t.Cell(16,3).Range.Text = "85-25="
$t.Cell(16,4).Range.Text = "3+52="
$t.Cell(16,5).Range.Text = "71-45="
$t.Cell(17,1).Range.Text = "36+18="
$t.Cell(17,2).Range.Text = "10+62="
$t.Cell(17,3).Range.Text = "48-19="
$t.Cell(17,4).Range.Text = "99-79="
$t.Cell(17,5).Range.Text = "99-11="
$t.Cell(18,1).Range.Text = "79-12="
$t.Cell(18,2).Range.Text = "68-27="
$t.Cell(18,3).Range.Text = "83-41="
$t.Cell(18,4).Range.Text = "2+20="
$t.Cell(18,5).Range.Text = "8+70="
$t.Cell(19,1).Range.Text = "71-13="
$t.Cell(19,2).Range.Text = "67-45="
$t.Cell(19,3).Range.Text = "60-34="
$t.Cell(19,4).Range.Text = "81+1="
$t.Cell(19,5).Range.Text = "93-73="
$t.Cell(20,1).Range.Text = "52-16="
$t.Cell(20,2).Range.Text = "36+20="
$t.Cell(20,3).Range.Text = "59+33="
$t.Cell(20,4).Range.Text = "57+26="
$t.Cell(20,5).Range.Text = "36-26="
